$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5273163333333333
$ws.Range("H2").Value = 1.581949
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 3.948345378425999
$ws.Range("R2").Value = 35.53510840583399
$ws.Range("S2").Value = 0.1384395179233961
$ws.Range("T2").Value = 0.1384395179233961

# Row 3
$ws.Range("G3").Value = 0.5273163333333333
$ws.Range("H3").Value = 1.581949
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 16.873821041724
$ws.Range("R3").Value = 151.864389375516
$ws.Range("S3").Value = 0.5916411627275552
$ws.Range("T3").Value = 0.5916411627275552

# Row 4
$ws.Range("G4").Value = 0.5273163333333333
$ws.Range("H4").Value = 1.581949
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 7.698197112253887
$ws.Range("R4").Value = 69.28377401028499
$ws.Range("S4").Value = 0.2699193193490487
$ws.Range("T4").Value = 0.2699193193490487
